# ft_printf testing in progress
# Applies the edits described in the commit:
#  - Feuil1!D7 gets a formula computing unsigned-int overflow wraparound
#  - modifs libft (sheet2) gets a new row documenting `ft_put_uint_fd`
#    and the existing "return int" note for ft_putstr_fd is expanded.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("modifs libft")

# --- Feuil1: D7 formula (unsigned int wraparound demonstration) ---
$ws1.Range("D7").Formula = "=-2147483648 + 4294967295+1"

# --- Feuil1: selection cosmetic change ---
$ws1.Range("B13:D18").Select()

# --- modifs libft: widen column B to fit the longer note ---
$ws2.Columns.Item(2).ColumnWidth = 33.140625

# --- modifs libft: update B3 (ft_putstr_fd's note) with the fuller text ---
$ws2.Range("B3").Value = "return int et write (null) si s == NULL"

# --- modifs libft: append new row describing the new ft_put_uint_fd func ---
$ws2.Range("A6").Value = "ft_put_uint_fd"
$ws2.Range("B6").Value = "new func"

# Give A6 a thin left/right border (no top/bottom) to match the other column-A cells' look
$ws2.Range("A6").Borders.Item(7).LineStyle = 1    # xlEdgeLeft / xlContinuous
$ws2.Range("A6").Borders.Item(7).Weight = 2       # xlThin
$ws2.Range("A6").Borders.Item(10).LineStyle = 1   # xlEdgeRight / xlContinuous
$ws2.Range("A6").Borders.Item(10).Weight = 2      # xlThin

# --- modifs libft: selection cosmetic change ---
$ws2.Range("A6").Select()

# --- Recalculate so cached formula values are written ---
$excel.Calculate()
